$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Reassign the batch owner (column C) for a handful of transaction table rows
$ws.Range("C24").Value = "Raghava/Divya"
$ws.Range("C25").Value = "Raghava"
$ws.Range("C30").Value = "Raghava"
$ws.Range("C33").Value = "Divya"

# New table entry added to the TRANSACTIONS list: Login
$ws.Range("A37").Value = "Login"

# Update the active selection to reflect where work is currently happening
$ws.Range("C31").Select()
